# Adds a new weekly price record for "Ajo" (Hortaliza) at
# Terminal La Palmera de La Serena. The new record is inserted as row 546,
# pushing all existing records from row 546 downward by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 546; everything below shifts down one row.
$ws.Rows.Item(546).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A546").Value2 = 8
$ws.Range("B546").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C546").Value2 = "Coquimbo"
$ws.Range("D546").Value2 = 45265
$ws.Range("E546").Value2 = 4
$ws.Range("F546").Value2 = 100112003
$ws.Range("G546").Value2 = "Ajo"
$ws.Range("H546").Value2 = "Chino"
$ws.Range("I546").Value2 = "Primera"
$ws.Range("J546").Value2 = 360
$ws.Range("K546").Value2 = 23000
$ws.Range("L546").Value2 = 24000
$ws.Range("M546").Value2 = 23500
$ws.Range("N546").Value2 = "`$/caja 10 kilos"
$ws.Range("O546").Value2 = "China"
$ws.Range("P546").Value2 = 2350
$ws.Range("Q546").Value2 = 10
$ws.Range("R546").Value2 = "Hortaliza"
